# Updated question with knowledge of net flow
#
# Original body (3 paragraphs):
#   P1 "Question:"
#   P2 "How does the type of donor of agricultural development aid designated
#       for nutrition influence the level of food security in countries of
#       Sub-Saharan Africa?"
#   P3 (empty, holds the _GoBack bookmark only)
#
# Target body (6 paragraphs):
#   P1 "Question:"                                        (unchanged)
#   P2 "How does the type of donor ... Africa?"            (unchanged)
#   P3 (new, empty)
#   P4 (new) "--Knowing it is net flows\u2014"
#   P5 (was P3) "How does the source of positive flows ... Africa?"
#       + the original bookmarkStart/bookmarkEnd
#       + a trailing run containing a single space
#   P6 (new, empty)

$d = $word.ActiveDocument

# --- Insert a new blank paragraph right after paragraph 2 ---------------
$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertParagraphAfter()

# --- Insert the "net flows" aside paragraph after that new blank one ----
$blank = $d.Paragraphs.Item(3)
$blank.Range.InsertParagraphAfter()
$aside = $d.Paragraphs.Item(4)
$aside.Range.Text = "--Knowing it is net flows" + [char]0x2014

# --- Flesh out the paragraph that still holds the _GoBack bookmark ------
# Insert the trailing space immediately after the bookmark FIRST (while the
# bookmark still sits at the very start of the otherwise-empty paragraph)
# so it lands in its own run, after bookmarkEnd.
$bm = $d.Bookmarks.Item("_GoBack")
$afterBookmark = $d.Range($bm.End, $bm.End)
$afterBookmark.InsertAfter(" ")

# Now insert the question text before the (still collapsed) bookmark so it
# becomes its own run ahead of bookmarkStart.
$bm = $d.Bookmarks.Item("_GoBack")
$beforeBookmark = $d.Range($bm.Start, $bm.Start)
$beforeBookmark.InsertBefore("How does the source of positive flows of agricultural development aid designated for nutrition influence the level of food security in countries of Sub-Saharan Africa?")

# --- Trailing blank paragraph at the very end ----------------------------
$questionPara = $d.Paragraphs.Item(5)
$questionPara.Range.InsertParagraphAfter()
